# Generate Report for handback
# Updates the status of the file "878e5e59-2378-49c6-bcb0-0a60edbf870d.md"
# from "Ready for handoff" to "Handed back: in sync with en-US" for both
# the zh-cn and de-de localization targets, along with the handback
# timestamps, and reflects that on the Overview sheet as well.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# Overview sheet: update the zh-cn / de-de status columns for row 3
# (878e5e59-2378-49c6-bcb0-0a60edbf870d.md)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $handedBack
$wsOverview.Range("C3").Value = $handedBack

# zh-cn sheet: update Status (B3) and Latest Handback DateTime (G3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $handedBack
$wsZhCn.Range("G3").Value = "2016-01-11 13:37:43"

# de-de sheet: update Status (B3) and Latest Handback DateTime (G3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $handedBack
$wsDeDe.Range("G3").Value = "2016-01-11 13:38:13"
